$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Date updates
#    "Entrega final: 01/09" -> "16/02"  (there are two "01/09" dates
#    in the document and both become "16/02", so Replace All is safe)
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute("01/09", $false, $false, $false, $false, $false, `
                                 $true, 1, $false, "16/02", 2)

# "28/07" -> "12/01" (single occurrence)
$null = $d.Content.Find.Execute("28/07", $false, $false, $false, $false, $false, `
                                 $true, 1, $false, "12/01", 2)

# ------------------------------------------------------------------
# 2) "grupos de 3 pessoas" -> "grupos de 3 ou 4 pessoas"
#    Insert " ou 4" right after the bold "3" run.
# ------------------------------------------------------------------
$r2 = $d.Content
$null = $r2.Find.Execute("3 pessoas")
$insertAt = $d.Range($r2.Start + 1, $r2.Start + 1)
$insertAt.InsertAfter(" ou 4")

# ------------------------------------------------------------------
# 3) "(os grupos deverão ..." -> "(Os grupos deverão ..."
#    Capitalize the "o" right after the opening parenthesis.
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute("(os grupos deverão ser os mesmos definidos no Trabalho 1)", `
                                 $true, $false, $false, $false, $false, $true, 1, $false, `
                                 "(Os grupos deverão ser os mesmos definidos no Trabalho 1)", 2)

# ------------------------------------------------------------------
# 4) Move the "_GoBack" bookmark from the end of the tab-filled
#    paragraph to the middle of the word "preferência" (splitting it
#    into "preferênci" + bookmark + "a"). Re-adding a bookmark with
#    the same name moves it, so the stray one at the old location is
#    removed automatically.
# ------------------------------------------------------------------
$r3 = $d.Content
$null = $r3.Find.Execute("preferência")
$splitPoint = $r3.End - 1
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
